$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column S (19th column) width 20 -> 21 (raw OOXML width units).
# COM ColumnWidth = raw_width - 5/6, so raw 21 => ColumnWidth 20.1666666666667
$ws.Columns.Item(19).ColumnWidth = 20.1666666666667

# Update recalculated values (rows 2-12) for columns P, Q, S, T, U, V, X, Z
$ws.Range("P2").Value = 0.4723303633292488
$ws.Range("Q2").Value = 140.2035592597601
$ws.Range("S2").Value = 0.08842508596596292
$ws.Range("T2").Value = 0.5747214936203638
$ws.Range("U2").Value = 0.07373597010466124
$ws.Range("V2").Value = 0.4792491452980848
$ws.Range("X2").Value = 0.2404719358073799
$ws.Range("Z2").Value = 0.474932499321271

$ws.Range("P3").Value = 0.4205686640219093
$ws.Range("Q3").Value = 117.8000861903223
$ws.Range("S3").Value = 0.07456781122585729
$ws.Range("T3").Value = 0.4573291322971118
$ws.Range("U3").Value = 0.06111404150919723
$ws.Range("V3").Value = 0.3748163063270806
$ws.Range("X3").Value = 0.212528263659296
$ws.Range("Z3").Value = 0.4126602370625386

$ws.Range("P4").Value = 0.5718412463218272
$ws.Range("Q4").Value = 178.4320803110084
$ws.Range("S4").Value = 0.1215599161033057
$ws.Range("T4").Value = 0.9099502856655908
$ws.Range("U4").Value = 0.1022438583907255
$ws.Range("V4").Value = 0.7653577851364012
$ws.Range("X4").Value = 0.2595318066581477
$ws.Range("Z4").Value = 0.5185024253006697

$ws.Range("P5").Value = 0.5581335994204907
$ws.Range("Q5").Value = 167.8992861418995
$ws.Range("S5").Value = 0.08827337731641663
$ws.Range("T5").Value = 0.6689697604315441
$ws.Range("U5").Value = 0.07201059816602934
$ws.Range("V5").Value = 0.5457241364062082
$ws.Range("X5").Value = 0.2394003956200058
$ws.Range("Z5").Value = 0.4526383053569616

$ws.Range("P6").Value = 0.7001746024968
$ws.Range("Q6").Value = 228.0927372156346
$ws.Range("S6").Value = 0.09483997034262538
$ws.Range("T6").Value = 0.8326105656603114
$ws.Range("U6").Value = 0.06698611160786123
$ws.Range("V6").Value = 0.5880784660277253
$ws.Range("X6").Value = 0.2767665980761518
$ws.Range("Z6").Value = 0.4202146271806473

$ws.Range("P7").Value = 1.019914157386467
$ws.Range("Q7").Value = 387.8307114221698
$ws.Range("S7").Value = 0.08019680609824066
$ws.Range("T7").Value = 0.8786003547184101
$ws.Range("U7").Value = 0.03213493902083905
$ws.Range("V7").Value = 0.3520560256224849
$ws.Range("X7").Value = 0.3662718240975267
$ws.Range("Z7").Value = 0.4204215952509643

$ws.Range("P8").Value = 1.07104920097898
$ws.Range("Q8").Value = 397.6694718519936
$ws.Range("S8").Value = 0.1177110321319092
$ws.Range("T8").Value = 1.386957062992251
$ws.Range("U8").Value = 0.09652239426565881
$ws.Range("V8").Value = 1.137297108342898
$ws.Range("X8").Value = 0.3461776989580811
$ws.Range("Z8").Value = 0.6402659933017059

$ws.Range("P9").Value = 1.361273927133048
$ws.Range("Q9").Value = 574.953156589773
$ws.Range("S9").Value = 0.128554978317513
$ws.Range("T9").Value = 1.692375571059443
$ws.Range("U9").Value = 0.09075829707541901
$ws.Range("V9").Value = 1.194797174342258
$ws.Range("X9").Value = 0.445119024444915
$ws.Range("Z9").Value = 0.6555079664338384

$ws.Range("P10").Value = 1.193679569246113
$ws.Range("Q10").Value = 459.1098458480349
$ws.Range("S10").Value = 0.1220241098281319
$ws.Range("T10").Value = 1.546870014982985
$ws.Range("U10").Value = 0.08756530346934376
$ws.Range("V10").Value = 1.110044092764901
$ws.Range("X10").Value = 0.3682598322302786
$ws.Range("Z10").Value = 0.5582969920311238

$ws.Range("P11").Value = 1.477936038466737
$ws.Range("Q11").Value = 605.0392506059865
$ws.Range("S11").Value = 0.3304522886183035
$ws.Range("T11").Value = 4.872880388587285
$ws.Range("U11").Value = 0.3197505609100884
$ws.Range("V11").Value = 4.715071709787069
$ws.Range("X11").Value = 0.4138235051315636
$ws.Range("Z11").Value = 1.723688684572985

$ws.Range("P12").Value = 1.643165843958127
$ws.Range("Q12").Value = 692.8114952454029
$ws.Range("S12").Value = 0.6385336411494568
$ws.Range("T12").Value = 10.16437674197399
$ws.Range("U12").Value = 0.6309002809125057
$ws.Range("V12").Value = 10.04286654383328
$ws.Range("X12").Value = 0.437297857237506
$ws.Range("Z12").Value = 2.97888207228676
